# Fixed import of PI and Investment instrument uniqueness
# Remove the unused "Folio No" column (column G) from the worksheet,
# shifting every column to its right one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Folio No" column (G). EntireColumn.Delete shifts all
# subsequent columns left, which matches the diff (H->G, I->H, J->I, K->J, L->K)
$ws.Range("G1").EntireColumn.Delete()
